$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''26.100.13'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.91%  '

# Row 3
$ws.Range("D3").Value = '''1.651.97'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.08%  '

# Row 4
$ws.Range("E4").Value = '  -0.56%  '

# Row 5
$ws.Range("D5").Value = '''218.67'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.83%  '

# Row 6
$ws.Range("D6").Value = '''0.5250'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.23%  '

# Row 7
$ws.Range("E7").Value = '  -0.50%  '

# Row 8
$ws.Range("D8").Value = '''0.2676'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.82%  '

# Row 9
$ws.Range("D9").Value = '''0.06369'
$ws.Range("D9").Style = "Normal"

# Row 10
$ws.Range("D10").Value = '''20.53'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.25%  '

# Row 11
$ws.Range("D11").Value = '''0.07689'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.10%  '

# Row 12
$ws.Range("D12").Value = '''4.597'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.20%  '

# Row 13
$ws.Range("D13").Value = '''1.664.66'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.14%  '

# Row 14
$ws.Range("D14").Value = '''1.879.07'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.00%  '

# Row 15
$ws.Range("D15").Value = '''0.5622'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.10%  '

# Row 16
$ws.Range("D16").Value = '''0.0₅8220'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.06%  '

# Row 17
$ws.Range("D17").Value = '''65.47'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.79%  '

# Row 18
$ws.Range("D18").Value = '''26.093.93'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.92%  '

# Row 20
$ws.Range("D20").Value = '''4.682'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.99%  '

# Row 21
$ws.Range("B21").Value = 'Avalanche'
$ws.Range("C21").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D21").Value = '''10.35'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.43%  '

# Row 22
$ws.Range("B22").Value = 'BitcoinCash'
$ws.Range("C22").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D22").Value = '''191.07'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -5.53%  '

# Row 23
$ws.Range("D23").Value = '''5.972'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.62%  '

# Row 24
$ws.Range("D24").Value = '''1.004'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.56%  '

# Row 25
$ws.Range("D25").Value = '''146.12'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.34%  '

# Row 26
$ws.Range("D26").Value = '''0.1202'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.23%  '

# Row 27
$ws.Range("D27").Value = '''7.240'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.21%  '

# Row 28
$ws.Range("D28").Value = '''15.95'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.77%  '

# Row 29
$ws.Range("D29").Value = '''1.498'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.56%  '

# Row 30
$ws.Range("D30").Value = '''0.05637'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.39%  '

# Row 31
$ws.Range("E31").Value = '  -1.34%  '

# Row 32
$ws.Range("D32").Value = '''3.493'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.22%  '

# Row 33
$ws.Range("D33").Value = '''3.383'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.46%  '

# Row 34
$ws.Range("D34").Value = '''1.578'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.72%  '

# Row 35
$ws.Range("D35").Value = '''2.793'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.38%  '

# Row 36
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").Value = '''0.9450'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.13%  '

# Row 37
$ws.Range("B37").Value = 'HuobiToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D37").Value = '''2.409'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.91%  '

# Row 38
$ws.Range("D38").Value = '''0.5772'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.59%  '

# Row 39
$ws.Range("D39").Value = '''0.01592'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.48%  '

# Row 40
$ws.Range("D40").Value = '''5.968'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.08%  '

# Row 41
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").Value = '''0.8443'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.68%  '

# Row 42
$ws.Range("B42").Value = 'PaxDollar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D42").Value = '''1.003'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.56%  '

# Row 43
$ws.Range("D43").Value = '''1.021.80'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -5.37%  '

# Row 44
$ws.Range("D44").Value = '''101.33'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.71%  '

# Row 45
$ws.Range("D45").Value = '''1.790.41'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.96%  '

# Row 46
$ws.Range("D46").Value = '''58.39'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.30%  '

# Row 47
$ws.Range("B47").Value = 'Frax'
$ws.Range("C47").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D47").Value = '''1.004'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.15%  '

# Row 48
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").Value = '''0.05342'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.74%  '

# Row 49
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '''8.037'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.58%  '

# Row 50
$ws.Range("B50").Value = 'Mantle'
$ws.Range("C50").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D50").Value = '''0.4341'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.68%  '

# Row 51
$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D51").Value = '''0.09744'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.40%  '
